$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (columns M-T)
$ws.Range("M2").Value = 0.1557005
$ws.Range("N2").Value = 0.311401
$ws.Range("O2").Value = 0.02102398211576467
$ws.Range("P2").Value = 0.01500040222529337
$ws.Range("Q2").Value = 0.089455853869
$ws.Range("R2").Value = 0.536735123214
$ws.Range("S2").Value = 0.02102398211576467
$ws.Range("T2").Value = 0.01500040222529337

# Update row 3 (columns O, P, S, T)
$ws.Range("O3").Value = 0.7998659708565604
$ws.Range("P3").Value = 0.8560446272575798
$ws.Range("S3").Value = 0.7998659708565604
$ws.Range("T3").Value = 0.8560446272575798

# Update row 4 (columns K-T)
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.024117
$ws.Range("N4").Value = 0.072351
$ws.Range("O4").Value = 0.003256478795417461
$ws.Range("P4").Value = 0.003485197868350457
$ws.Range("Q4").Value = 0.013856132946
$ws.Range("R4").Value = 0.124705196514
$ws.Range("S4").Value = 0.003256478795417461
$ws.Range("T4").Value = 0.003485197868350457

# Update row 5 (columns M-T)
$ws.Range("M5").Value = 1.3023455
$ws.Range("N5").Value = 2.604691
$ws.Range("O5").Value = 0.1758535682322574
$ws.Range("P5").Value = 0.1254697726487764
$ws.Range("Q5").Value = 0.7482469788789999
$ws.Range("R5").Value = 4.489481873273999
$ws.Range("S5").Value = 0.1758535682322574
$ws.Range("T5").Value = 0.1254697726487764

# Delete rows 6 and 7 (entire rows)
$ws.Rows("6:7").Delete()
